$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 2 and clear the formatting it inherited from the insert
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()
$ws.Range("A2").Value2 = -0.01545035839080799
$ws.Range("B2").Value2 = 0.2315296996384856
$ws.Range("C2").Value2 = -0.5798605158925052

# Append new rows 23-31 with data
$ws.Range("A23").Value2 = -1.618811368942268
$ws.Range("B23").Value2 = -0.6589505374431646
$ws.Range("C23").Value2 = 0.2403407692909387
$ws.Range("A24").Value2 = -0.08240008354186718
$ws.Range("B24").Value2 = 0.4783504903316493
$ws.Range("C24").Value2 = -3.809414207935333
$ws.Range("A25").Value2 = -3.95973014831543
$ws.Range("B25").Value2 = 0.9762580394744872
$ws.Range("C25").Value2 = -4.069071769714356
$ws.Range("A26").Value2 = -3.901577949523926
$ws.Range("B26").Value2 = 1.771272063255311
$ws.Range("C26").Value2 = -2.484678864479062
$ws.Range("A27").Value2 = -1.957046031951897
$ws.Range("B27").Value2 = -0.6577051877975557
$ws.Range("C27").Value2 = -7.9572014808655
$ws.Range("A28").Value2 = 0.0412573218345611
$ws.Range("B28").Value2 = -3.739429324865336
$ws.Range("C28").Value2 = -3.584903955459609
$ws.Range("A29").Value2 = -2.076164960861222
$ws.Range("B29").Value2 = -2.966795355081547
$ws.Range("C29").Value2 = -0.5832877159118695
$ws.Range("A30").Value2 = -5.426012933254244
$ws.Range("B30").Value2 = -0.3007338047027568
$ws.Range("C30").Value2 = 0.185311913490301
$ws.Range("A31").Value2 = -6.20224690437317
$ws.Range("B31").Value2 = 0.8901370018720657
$ws.Range("C31").Value2 = 1.72858691215514
